$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STM32L031_TSSOP20_Dev_Board_rev")

# Row 21 holds the Crystal (Y1) part. Both the "Comment" (A) and
# "DesignItemId" (E) columns reference the same part number string,
# which is being updated from the old FC-135 crystal to Q13FC1350000200.
$oldValue = "FC-135_32.7680KA-AG5"
$newValue = "Q13FC1350000200"

for ($row = 1; $row -le $ws.UsedRange.Rows.Count; $row++) {
    foreach ($col in @("A", "E")) {
        $cell = $ws.Range("$col$row")
        if ($cell.Value2 -eq $oldValue) {
            # Prefix with an apostrophe so Excel keeps treating the cell as
            # explicit text (preserves the original quote-prefixed style)
            # while the stored string itself stays clean.
            $cell.Value = "'" + $newValue
        }
    }
}
